$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.933979272842407
$ws.Range("B1").Value = 2.653825998306274
$ws.Range("C1").Value = 2.87424635887146
$ws.Range("D1").Value = 3.360890865325928
$ws.Range("E1").Value = 1.038919687271118
